$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, shifting existing rows 139:150 down to 140:151
$ws.Rows("139:139").Insert(4)

# Populate the newly inserted row 139 with the new weekly data point
$ws.Cells.Item(139, 1).Value = 4
$ws.Cells.Item(139, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(139, 3).Value = 'Los Lagos'
$ws.Cells.Item(139, 4).Value = 45142
$ws.Cells.Item(139, 5).Value = 10
$ws.Cells.Item(139, 6).Value = 100112026
$ws.Cells.Item(139, 7).Value = 'Haba'
$ws.Cells.Item(139, 8).Value = 'Sin especificar'
$ws.Cells.Item(139, 9).Value = 'Primera'
$ws.Cells.Item(139, 10).Value = 100
$ws.Cells.Item(139, 11).Value = 18000
$ws.Cells.Item(139, 12).Value = 20000
$ws.Cells.Item(139, 13).Value = 19000
$ws.Cells.Item(139, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(139, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(139, 16).Value = 760
$ws.Cells.Item(139, 17).Value = 25
$ws.Cells.Item(139, 18).Value = 'Hortaliza'
